$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.681.20'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '1.638.69'

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.525'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.74%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.10'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.25%  '

$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("E10").Value = '  -0.12%  '

$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D12").Value = '1.871.06'
$ws.Range("E12").Value = '  -0.53%  '

$ws.Range("D13").Value = '1.637.59'
$ws.Range("E13").Value = '  -1.31%  '

$ws.Range("E14").Value = '  +0.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.560'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.28%  '

$ws.Range("D17").Value = '27.666.96'
$ws.Range("E17").Value = '  +0.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.76%  '

$ws.Range("D20").Value = '0.0₃0723'
$ws.Range("E20").Value = '  -0.58%  '

$ws.Range("E22").Value = '  -0.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.98'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.68%  '

$ws.Range("E26").Value = '  -1.12%  '

$ws.Range("E27").Value = '  -1.66%  '

$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("E32").Value = '  -0.30%  '

$ws.Range("D33").Value = '1.458.12'
$ws.Range("E33").Value = '  +2.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.13'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.11%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("E37").Value = '  -0.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.877'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.40%  '

$ws.Range("E39").Value = '  +0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.898'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.74'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.96%  '

$ws.Range("E42").Value = '  -0.61%  '

$ws.Range("E43").Value = '  +0.06%  '

$ws.Range("E44").Value = '  +1.05%  '

$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.73%  '

$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.780.91'
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.74'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.87%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.65%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0107'
$ws.Range("E49").Value = '  -0.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0992'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.29%  '
